# Weekly refresh: insert the newest week's record for Zapallo / Paine at
# the top of the historical block (row 281) and push the existing rows
# down by one. The previously-last row (342) is pushed to 343 and the
# row that used to be at 341 becomes the new row 342.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 281 - shifts rows 281:342 down to 282:343
$ws.Rows.Item(281).Insert()

# Populate the newly inserted row with the new week's data
$ws.Range("A281").Value = 4
$ws.Range("B281").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C281").Value = "Los Lagos"
$ws.Range("D281").Value = 44782
$ws.Range("E281").Value = 10
$ws.Range("F281").Value = 100112045
$ws.Range("G281").Value = "Zapallo"
$ws.Range("H281").Value = "Paine"
$ws.Range("I281").Value = "1a (guarda)"
$ws.Range("J281").Value = 1200
$ws.Range("K281").Value = 550
$ws.Range("L281").Value = 650
$ws.Range("M281").Value = 600
$ws.Range("N281").Value = "$/kilo (volumen en unidades)"
$ws.Range("O281").Value = "Región de O'Higgins"
$ws.Range("P281").Value = 600
$ws.Range("Q281").Value = 1
$ws.Range("R281").Value = "Hortaliza"
